$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.266.91"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "3.501.52"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'586.45"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("D6").Value = "'134.32"
$ws.Range("E6").Value = "  +2.17%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.486"
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("E9").Value = "  +2.27%  "
$ws.Range("D10").Value = "'7.28"
$ws.Range("E10").Value = "  +2.21%  "
$ws.Range("D11").Value = "'0.386"
$ws.Range("E11").Value = "  +1.80%  "
$ws.Range("D12").Value = "4.098.50"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("E14").Value = "  +3.31%  "
$ws.Range("D15").Value = "3.501.32"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "'26.15"
$ws.Range("E16").Value = "  -4.42%  "
$ws.Range("D17").Value = "64.311.86"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").Value = "'13.78"
$ws.Range("E20").Value = "  -3.40%  "
$ws.Range("D21").Value = "'393.22"
$ws.Range("E21").Value = "  +2.60%  "
$ws.Range("D22").Value = "'0.572"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").Value = "3.641.73"
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("D24").Value = "'74.21"
$ws.Range("E24").Value = "  +2.08%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "'5.72"
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").Value = "'0.0000115"
$ws.Range("E27").Value = "  +2.38%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'7.53"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "'1.51"
$ws.Range("E30").Value = "  -3.46%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'8.29"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'2.23"
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("B33").Value = "RenzoRestakedETH"
$ws.Range("C33").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D33").Value = "3.523.38"
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.149"
$ws.Range("E35").Value = "  +3.71%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "'23.48"
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'5.21"
$ws.Range("E37").Value = "  -1.83%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'1.57"
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "'6.92"
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "'161.96"
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "'0.0782"
$ws.Range("E41").Value = "  -1.71%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "'0.806"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'25.46"
$ws.Range("E44").Value = "  -3.23%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "'4.41"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("E46").Value = "  +1.88%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "'1.17"
$ws.Range("E47").Value = "  -2.71%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.471.23"
$ws.Range("E48").Value = "  +2.47%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'6.80"
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").Value = "'0.895"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0262"
$ws.Range("E51").Value = "  -0.74%  "
